# Adjust Investment Summary table column widths for better formatting
# (and clear out the now-stale template placeholder copy/fill in those
#  comparison / financial-impact / risk tables)

$p = $ppt.ActivePresentation

function Clear-TableCells($tbl, $rowCount, $colCount) {
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $tbl.Cell($r, $c)
            $cell.Shape.TextFrame.TextRange.Text = ""
        }
    }
}

# EMU -> point helper (1 pt = 12700 EMU)
$emuPerPt = 12700.0

# --- Slide 2: "Traditional Approach" / "Our Solution" comparison table ---
$s2 = $p.Slides.Item(2)
$tbl2 = $s2.Shapes.Item(3).Table
Clear-TableCells $tbl2 4 2
$tbl2.Columns.Item(2).Width = 4355467 / $emuPerPt

# --- Slide 3: "Metric" / "Value" financial impact table ---
$s3 = $p.Slides.Item(3)
$tbl3 = $s3.Shapes.Item(3).Table
Clear-TableCells $tbl3 6 2
$tbl3.Columns.Item(2).Width = 4355467 / $emuPerPt

# --- Slide 4: "Risk" / "Mitigation Strategy" / "Success Probability" table ---
$s4 = $p.Slides.Item(4)
$tbl4 = $s4.Shapes.Item(3).Table
Clear-TableCells $tbl4 4 3
$tbl4.Columns.Item(3).Width = 2903645 / $emuPerPt

Write-Host "Investment Summary table column widths adjusted."
